# Refresh the "cryptos" price list (GitHub Actions-style scheduled update).
# Price/Volume cells hold text (not numbers), and a couple of rows got
# re-ranked (their Coin/Link/Price/Volume swapped with a neighboring row).
# For Price cells whose new text looks numeric, force text formatting first
# (then clear the format again) so Excel stores the exact literal string
# instead of silently coercing it to a float (e.g. "33.01" -> 33.00999...,
# "1.00" -> 1, "0.0662" -> 6.62E-2) and without leaving a residual style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.433.30"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.787.75"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.01"
$ws.Range("D8").ClearFormats()
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0662"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "2.045.81"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.00"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.92%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.788.88"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").Value = "34.392.20"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.19"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "255.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").Value = "0.0₃0746"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -4.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.46"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("D35").Value = "1.455.14"
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.632"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0188"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.892"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0507"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.89"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.944.52"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.09"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.29"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.58%  "
